$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: 8,-1)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: 7,-7)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: -5,-6)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: 8,8)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: -10,-1)"
$ws.Range("F1").Value = "(305251175, Or  Leder: 6,6)"
$ws.Range("G1").Value = "(308051846, Eyal  Sofer: 3,-3)"

$ws.Range("A3").Value = "cost: 867.1581816446853"
$ws.Range("A4").Value = "time: 105.26977270558567"
